$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 2504.364414539994
$ws.Range("B1").Value = 1675.5675376436604
$ws.Range("C1").Value = 1659.1065227383422
$ws.Range("A2").Value = 2263.2649011438675
$ws.Range("B2").Value = 1537.5762282943956
$ws.Range("C2").Value = 1393.3778585627351
$ws.Range("A3").Value = 2560.3680356085983
$ws.Range("B3").Value = 1740.3097987834262
$ws.Range("C3").Value = 1583.9747961737785
$ws.Range("A4").Value = 2495.9915444752205
$ws.Range("B4").Value = 1913.4487934162717
$ws.Range("C4").Value = 1951.9452111165515
$ws.Range("A5").Value = 2528.80384513526
$ws.Range("B5").Value = 1753.0003849402135
$ws.Range("C5").Value = 1754.546258168725
$ws.Range("A6").Value = 2480.0509760227787
$ws.Range("B6").Value = 1850.4755988030256
$ws.Range("C6").Value = 1931.3762670675633
$ws.Range("A7").Value = 2397.2088802092485
$ws.Range("B7").Value = 1848.4300945850457
$ws.Range("C7").Value = 1675.6680969497306
$ws.Range("A8").Value = 2467.9969987293316
$ws.Range("B8").Value = 1930.0304210380705
$ws.Range("C8").Value = 1785.2202694457196
$ws.Range("A9").Value = 2643.9044689181596
$ws.Range("B9").Value = 1947.2351251927282
$ws.Range("C9").Value = 1640.0636790523515
$ws.Range("A10").Value = 2394.9108734696329
$ws.Range("B10").Value = 1511.42538690984
$ws.Range("C10").Value = 1431.2601301716491
$ws.Range("A11").Value = 2163.0420361799524
$ws.Range("B11").Value = 1575.4238379800943
$ws.Range("C11").Value = 1407.2614046763103
$ws.Range("A12").Value = 2770.626449954108
$ws.Range("B12").Value = 2193.8781744486687
$ws.Range("C12").Value = 1875.1042057517848
$ws.Range("A13").Value = 2528.2387490187198
$ws.Range("B13").Value = 1947.5466896032067
$ws.Range("C13").Value = 1749.6330418674784
$ws.Range("A14").Value = 2609.7077573395409
$ws.Range("B14").Value = 2018.5468677722945
$ws.Range("C14").Value = 1759.5597087708768
$ws.Range("A15").Value = 2503.2297760451011
$ws.Range("B15").Value = 2051.416880995745
$ws.Range("C15").Value = 1833.0290098974242
$ws.Range("A16").Value = 2592.3768062015702
$ws.Range("B16").Value = 1797.2456078784303
$ws.Range("C16").Value = 1568.6528356072297
$ws.Range("A17").Value = 2356.8601401089932
$ws.Range("B17").Value = 1814.1836892896029
$ws.Range("C17").Value = 1709.7092249325049
$ws.Range("A18").Value = 2615.491122412915
$ws.Range("B18").Value = 2163.1508444250139
$ws.Range("C18").Value = 2077.7982519417487
$ws.Range("A19").Value = 2011.2261617487509
$ws.Range("B19").Value = 2033.9335197896344
$ws.Range("C19").Value = 1981.0698214975762
$ws.Range("A20").Value = 2576.2845513363582
$ws.Range("B20").Value = 1981.3717182002993
$ws.Range("C20").Value = 1870.0558660754841
$ws.Range("A21").Value = 2713.0802018334111
$ws.Range("B21").Value = 2020.8570458484212
$ws.Range("C21").Value = 1980.6742375270492
$ws.Range("A22").Value = 2569.5114395448936
$ws.Range("B22").Value = 1982.4411079456511
$ws.Range("C22").Value = 1731.9360952080569
